$wb = $excel.ActiveWorkbook

# Guns sheet keeps its data; only the selection/active-cell changes.
$gunsSheet = $wb.Worksheets.Item("Guns")
$gunsSheet.Range("B35").Select() | Out-Null

# Add the new "Slots" sheet right after "Guns".
$ws = $wb.Worksheets.Add($null, $gunsSheet)
$ws.Name = "Slots"

# Populate column A: header "name" followed by the slot types (alphabetical).
$values = @(
    "name",
    "Barrel",
    "Bipod",
    "Charge",
    "Dust Cover",
    "Gas Mod",
    "Handguard",
    "Magazine",
    "Mount",
    "Muzzle",
    "Optic",
    "Pistol Grip",
    "Receiver",
    "Sight",
    "Stock",
    "Tactical Grip",
    "Tactical Mod"
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Header row bold, like the other sheet.
$ws.Range("A1").Font.Bold = $true

# Column A width ~21 characters.
$ws.Columns.Item(1).ColumnWidth = 20.1666666667

# Match default page setup used across the workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection on the new sheet.
$ws.Range("F5").Select() | Out-Null

Write-Output "done"
